$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.860.43"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.568.02"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'576.60"
$ws.Range("E5").Value = "  -3.27%  "
$ws.Range("D6").Value = "'187.59"
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "  -3.33%  "
$ws.Range("D8").Value = "3.564.57"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -3.16%  "
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "'55.75"
$ws.Range("E12").Value = "  -4.36%  "
$ws.Range("D13").Value = "'0.0000298"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "4.137.11"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "'19.88"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "3.562.25"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "69.707.36"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'12.52"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").Value = "'471.61"
$ws.Range("E22").Value = "  -5.39%  "
$ws.Range("D23").Value = "'19.11"
$ws.Range("E23").Value = "  +12.70%  "
$ws.Range("D24").Value = "'5.06"
$ws.Range("E24").Value = "  -8.77%  "
$ws.Range("D25").Value = "'4.33"
$ws.Range("E25").Value = "  -3.46%  "
$ws.Range("D26").Value = "'88.30"
$ws.Range("E26").Value = "  -3.51%  "
$ws.Range("D27").Value = "'3.04"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").Value = "'10.94"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("D29").Value = "'9.30"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "'32.03"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("D33").Value = "'12.04"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").Value = "'65.64"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'570.67"
$ws.Range("E35").Value = "  -7.71%  "
$ws.Range("D36").Value = "'38.51"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "0.0₃0798"
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D39").Value = "'0.395"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("D40").Value = "'0.139"
$ws.Range("E40").Value = "  -6.24%  "
$ws.Range("D41").Value = "'3.48"
$ws.Range("E41").Value = "  -6.02%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'2.85"
$ws.Range("E42").Value = "  +5.23%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.209.91"
$ws.Range("E43").Value = "  -3.78%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'3.11"
$ws.Range("E44").Value = "  +10.55%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'3.10"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").Value = "'0.0440"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "'9.47"
$ws.Range("E47").Value = "  +4.33%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "'3.14"
$ws.Range("E51").Value = "  -3.73%  "
